$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 289, shifting existing rows 289-300 down to 290-301.
$ws.Rows.Item(289).Insert()

# Populate the new row 289 with the new price-record data.
$ws.Cells.Item(289, 1).Value2 = 5
$ws.Cells.Item(289, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(289, 3).Value2 = "Maule"
$ws.Cells.Item(289, 4).Value2 = 45075
$ws.Cells.Item(289, 5).Value2 = 7
$ws.Cells.Item(289, 6).Value2 = 100112017
$ws.Cells.Item(289, 7).Value2 = "Apio"
$ws.Cells.Item(289, 8).Value2 = "Americana (o)"
$ws.Cells.Item(289, 9).Value2 = "Primera"
$ws.Cells.Item(289, 10).Value2 = 700
$ws.Cells.Item(289, 11).Value2 = 5000
$ws.Cells.Item(289, 12).Value2 = 5000
$ws.Cells.Item(289, 13).Value2 = 5000
$ws.Cells.Item(289, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(289, 15).Value2 = "Provincia del Elquí"
$ws.Cells.Item(289, 16).Value2 = 833
$ws.Cells.Item(289, 17).Value2 = 6
$ws.Cells.Item(289, 18).Value2 = "Hortaliza"
